# Scheduled market-data refresh: updates currentAveragePrice* / Leve price &
# profit columns (H:N) across the leve-profit worksheets to the latest
# Universalis snapshot values.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 53999.5
$ws.Range("J75").Value = 53999.5
$ws.Range("L75").Value = 53999.5
$ws.Range("N75").Value = -55871.5
$ws.Range("H76").Value = 3325
$ws.Range("I76").Value = 3295.652
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 3295.652
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -2980.652
$ws.Range("N76").Value = -4630
$ws.Range("H78").Value = 53999.5
$ws.Range("J78").Value = 53999.5
$ws.Range("L78").Value = 161998.5
$ws.Range("N78").Value = -171358.5
$ws.Range("H79").Value = 3325
$ws.Range("I79").Value = 3295.652
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 3295.652
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -2203.652
$ws.Range("N79").Value = -6184
$ws.Range("H111").Value = 988.6667
$ws.Range("J111").Value = 1116
$ws.Range("L111").Value = 3348
$ws.Range("N111").Value = -9482
$ws.Range("H129").Value = 1630
$ws.Range("I129").Value = 579.6
$ws.Range("J129").Value = 2438
$ws.Range("K129").Value = 1738.8
$ws.Range("L129").Value = 7314
$ws.Range("M129").Value = 3261.2
$ws.Range("N129").Value = -17314
$ws.Range("H131").Value = 2555
$ws.Range("I131").Value = 1808.3334
$ws.Range("J131").Value = 3035
$ws.Range("K131").Value = 5425.0002
$ws.Range("L131").Value = 9105
$ws.Range("M131").Value = -385.0002000000004
$ws.Range("N131").Value = -19185
$ws.Range("H132").Value = 259118.16
$ws.Range("I132").Value = 297011.8
$ws.Range("J132").Value = 1441.2
$ws.Range("K132").Value = 891035.3999999999
$ws.Range("L132").Value = 4323.6
$ws.Range("M132").Value = -888505.3999999999
$ws.Range("N132").Value = -9383.6
$ws.Range("H135").Value = 26318692
$ws.Range("I135").Value = 1333.7693
$ws.Range("J135").Value = 83339630
$ws.Range("K135").Value = 12003.9237
$ws.Range("L135").Value = 750056670
$ws.Range("M135").Value = -9468.923699999999
$ws.Range("N135").Value = -750061740

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3269914.8
$ws.Range("I2").Value = 2352.8572
$ws.Range("J2").Value = 14706382
$ws.Range("K2").Value = 2352.8572
$ws.Range("L2").Value = 14706382
$ws.Range("M2").Value = -2239.8572
$ws.Range("N2").Value = -14706608
$ws.Range("H45").Value = 964.9375
$ws.Range("I45").Value = 808.9
$ws.Range("J45").Value = 1225
$ws.Range("K45").Value = 808.9
$ws.Range("L45").Value = 1225
$ws.Range("M45").Value = -431.9
$ws.Range("N45").Value = -1979
$ws.Range("H110").Value = 1024.2
$ws.Range("I110").Value = 862.58826
$ws.Range("J110").Value = 1940
$ws.Range("K110").Value = 862.58826
$ws.Range("L110").Value = 1940
$ws.Range("M110").Value = 1182.41174
$ws.Range("N110").Value = -6030
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H116").Value = 3269914.8
$ws.Range("I116").Value = 2352.8572
$ws.Range("J116").Value = 14706382
$ws.Range("K116").Value = 2352.8572
$ws.Range("L116").Value = 14706382
$ws.Range("M116").Value = -58.85719999999992
$ws.Range("N116").Value = -14710970
$ws.Range("H118").Value = 30800
$ws.Range("J118").Value = 30800
$ws.Range("L118").Value = 30800
$ws.Range("N118").Value = -34114
$ws.Range("H119").Value = 29122.143
$ws.Range("J119").Value = 29122.143
$ws.Range("L119").Value = 29122.143
$ws.Range("N119").Value = -38798.143
$ws.Range("H123").Value = 49320
$ws.Range("J123").Value = 49320
$ws.Range("L123").Value = 49320
$ws.Range("N123").Value = -59120
$ws.Range("H125").Value = 42989.75
$ws.Range("J125").Value = 42989.75
$ws.Range("L125").Value = 42989.75
$ws.Range("N125").Value = -52829.75
$ws.Range("H132").Value = 2302344
$ws.Range("I132").Value = 2876930.2
$ws.Range("K132").Value = 8630790.600000001
$ws.Range("M132").Value = -8628260.600000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3269914.8
$ws.Range("I3").Value = 2352.8572
$ws.Range("J3").Value = 14706382
$ws.Range("K3").Value = 2352.8572
$ws.Range("L3").Value = 14706382
$ws.Range("M3").Value = -2238.8572
$ws.Range("N3").Value = -14706610
$ws.Range("H134").Value = 25933.355
$ws.Range("I134").Value = 27034.906
$ws.Range("K134").Value = 81104.71799999999
$ws.Range("M134").Value = -78569.71799999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1494.7
$ws.Range("I99").Value = 1243.375
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 1243.375
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = 254.625
$ws.Range("N99").Value = -5496
$ws.Range("H122").Value = 1572.56
$ws.Range("I122").Value = 1408.8
$ws.Range("J122").Value = 2227.6
$ws.Range("K122").Value = 4226.4
$ws.Range("L122").Value = 6682.799999999999
$ws.Range("M122").Value = -1776.4
$ws.Range("N122").Value = -11582.8
$ws.Range("H126").Value = 1494.7
$ws.Range("I126").Value = 1243.375
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 3730.125
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -1260.125
$ws.Range("N126").Value = -12440

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3251.838
$ws.Range("I113").Value = 531.8946999999999
$ws.Range("J113").Value = 6122.8887
$ws.Range("K113").Value = 1595.6841
$ws.Range("L113").Value = 18368.6661
$ws.Range("M113").Value = 574.3159000000001
$ws.Range("N113").Value = -22708.6661
$ws.Range("H131").Value = 1889464.4
$ws.Range("J131").Value = 2274741.8
$ws.Range("L131").Value = 6824225.399999999
$ws.Range("N131").Value = -6834305.399999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 4000
$ws.Range("I35").Value = 4000
$ws.Range("K35").Value = 4000
$ws.Range("M35").Value = -3702
$ws.Range("H102").Value = 1014.17645
$ws.Range("I102").Value = 745.3
$ws.Range("J102").Value = 1398.2858
$ws.Range("K102").Value = 745.3
$ws.Range("L102").Value = 1398.2858
$ws.Range("M102").Value = 876.7
$ws.Range("N102").Value = -4642.2858
$ws.Range("H113").Value = 1936.5
$ws.Range("I113").Value = 2035
$ws.Range("J113").Value = 1739.5
$ws.Range("K113").Value = 2035
$ws.Range("L113").Value = 1739.5
$ws.Range("M113").Value = 135
$ws.Range("N113").Value = -6079.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 18167.5
$ws.Range("J14").Value = 18167.5
$ws.Range("L14").Value = 18167.5
$ws.Range("N14").Value = -18503.5
$ws.Range("H107").Value = 645.05884
$ws.Range("I107").Value = 371.45456
$ws.Range("J107").Value = 1146.6666
$ws.Range("K107").Value = 1114.36368
$ws.Range("L107").Value = 3439.9998
$ws.Range("M107").Value = 805.6363200000001
$ws.Range("N107").Value = -7279.9998
$ws.Range("H116").Value = 31000
$ws.Range("J116").Value = 31000
$ws.Range("L116").Value = 31000
$ws.Range("N116").Value = -40178
$ws.Range("H126").Value = 3777.5557
$ws.Range("I126").Value = 4187.875
$ws.Range("J126").Value = 495
$ws.Range("K126").Value = 12563.625
$ws.Range("L126").Value = 1485
$ws.Range("M126").Value = -10093.625
$ws.Range("N126").Value = -6425
$ws.Range("H133").Value = 35000
$ws.Range("J133").Value = 35000
$ws.Range("L133").Value = 35000
$ws.Range("N133").Value = -45120

